# Update "想去人数" (F column) values as published at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 4..9, column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1876
$ws1.Range("F5").Value = 1568
$ws1.Range("F6").Value = 285
$ws1.Range("F7").Value = 73
$ws1.Range("F8").Value = 552
$ws1.Range("F9").Value = 131

# Sheet "全部类型" - rows 4..6 and 8..10, column F (row 7 is an unrelated concert entry)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1876
$ws4.Range("F5").Value = 1568
$ws4.Range("F6").Value = 285
$ws4.Range("F8").Value = 73
$ws4.Range("F9").Value = 552
$ws4.Range("F10").Value = 131
